$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Requisitos:" row (23) and its value row (24)
$ws.Rows("23:24").Delete()

# Update "Semestre ideal" value from "EP-10" to "EA-1,EP-10"
$ws.Range("B9:C9").Value = "EA-1,EP-10"
